# Duplicate the "Chocolate sweet roll..." / "Bear claw dessert..." block
# (with its blank-paragraph separators) and insert two copies of it
# right after the "Tiramisu wafer cupcake..." paragraph.

$d = $word.ActiveDocument

# Locate the anchor paragraph: the one ending "...powder toffee."
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Tiramisu wafer cupcake*") {
        $anchor = $p
        break
    }
}

$newParaTexts = @(
    "",
    "Chocolate sweet roll cake gummi bears muffin. Candy canes pie powder cake macaroon cake. Lollipop sugar plum sesame snaps toffee soufflé. Topping marshmallow jelly apple pie toffee pie jelly dessert. Shortbread donut pudding chocolate cake bonbon cupcake tootsie roll gummi bears. Marshmallow jelly chupa chups wafer jelly beans marzipan cotton candy cookie. Pudding tootsie roll topping caramels cheesecake gummies tootsie roll caramels biscuit. Jelly candy dessert macaroon bear claw sweet roll carrot cake jelly beans. Topping sesame snaps biscuit tart sugar plum wafer macaroon gummies. Donut gummies dragée tart jelly beans cake caramels lollipop oat cake. Pudding soufflé pie marzipan candy. Sesame snaps dessert cupcake oat cake cake lollipop topping. Wafer sesame snaps pudding donut soufflé halvah gummies. Chocolate bar fruitcake macaroon sesame snaps powder cheesecake jujubes.",
    "",
    "Bear claw dessert dragée candy canes muffin. Marshmallow pudding pastry soufflé gummi bears cookie. Powder gummies shortbread tiramisu marzipan biscuit gummies soufflé. Carrot cake chupa chups candy canes topping chocolate jelly lemon drops toffee. Fruitcake soufflé tiramisu macaroon muffin gummies cheesecake. Cupcake biscuit bear claw bear claw chocolate cheesecake. Cheesecake shortbread lollipop carrot cake liquorice toffee pudding sugar plum tiramisu. Marshmallow shortbread sweet cheesecake tiramisu shortbread cheesecake cake candy canes. Jelly shortbread sweet brownie lemon drops tart pudding pie. Apple pie dessert cheesecake donut caramels. Gummi bears fruitcake cake cheesecake cake wafer chupa chups tart. Cheesecake jujubes chocolate cake candy canes wafer.",
    "",
    "Chocolate sweet roll cake gummi bears muffin. Candy canes pie powder cake macaroon cake. Lollipop sugar plum sesame snaps toffee soufflé. Topping marshmallow jelly apple pie toffee pie jelly dessert. Shortbread donut pudding chocolate cake bonbon cupcake tootsie roll gummi bears. Marshmallow jelly chupa chups wafer jelly beans marzipan cotton candy cookie. Pudding tootsie roll topping caramels cheesecake gummies tootsie roll caramels biscuit. Jelly candy dessert macaroon bear claw sweet roll carrot cake jelly beans. Topping sesame snaps biscuit tart sugar plum wafer macaroon gummies. Donut gummies dragée tart jelly beans cake caramels lollipop oat cake. Pudding soufflé pie marzipan candy. Sesame snaps dessert cupcake oat cake cake lollipop topping. Wafer sesame snaps pudding donut soufflé halvah gummies. Chocolate bar fruitcake macaroon sesame snaps powder cheesecake jujubes.",
    "",
    "Bear claw dessert dragée candy canes muffin. Marshmallow pudding pastry soufflé gummi bears cookie. Powder gummies shortbread tiramisu marzipan biscuit gummies soufflé. Carrot cake chupa chups candy canes topping chocolate jelly lemon drops toffee. Fruitcake soufflé tiramisu macaroon muffin gummies cheesecake. Cupcake biscuit bear claw bear claw chocolate cheesecake. Cheesecake shortbread lollipop carrot cake liquorice toffee pudding sugar plum tiramisu. Marshmallow shortbread sweet cheesecake tiramisu shortbread cheesecake cake candy canes. Jelly shortbread sweet brownie lemon drops tart pudding pie. Apple pie dessert cheesecake donut caramels. Gummi bears fruitcake cake cheesecake cake wafer chupa chups tart. Cheesecake jujubes chocolate cake candy canes wafer."
)

$cur = $anchor
foreach ($t in $newParaTexts) {
    $cur.Range.InsertParagraphAfter()
    $cur = $cur.Next()
    if ($t -ne "") {
        $cur.Range.Text = $t
    }
}
